# This script applies the updated FIM (Fiscal Impact Model) projection values
# for the "old order" run, per the commit:
#   "Run the FIM with the old order
#    Unfortunately, the projections do not match. Will have to investigate."
#
# It updates a block of quarterly projection cells (rows 219-260, columns
# D..AG) on the active worksheet to the new recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each line: Row,ColumnLetter,NewValue
$updates = @"
219,D,0.00814832810809674
219,E,0.00991870914463755
219,F,0.00991870914463755
219,G,0.00991870914463755
219,H,0.00784098294506985
219,J,28658.5
219,K,1863.8
219,L,433.054109
219,N,3137.4
219,O,4439
219,P,2328.3
219,T,2070.499
219,X,0.300000000000001
219,Z,92.635
219,AF,1606.2
219,AG,319.7
220,D,0.00631637071565883
220,E,0.00563330429941056
220,F,0.00563330429941056
220,G,0.00563330429941056
220,J,28968.9
220,K,1862.79273174517
220,L,437.682352782907
220,N,3166.71347588412
220,O,4488.43810791123
220,P,2365.66523074806
220,T,2077.43233333333
220,Z,92.635
220,AF,1632.88132683418
220,AG,322.260664845667
221,D,0.00580887002259489
221,E,0.00534045807802097
221,F,0.00534045807802097
221,G,0.00534045807802097
221,J,29256.5
221,K,1881.18376071914
221,L,442.603153049303
221,N,3195.65392542925
221,O,4537.15803906129
221,P,2389.58107566087
221,T,2085.89366666667
221,Z,99.734
221,AF,1648.41368264403
221,AG,323.696248231545
222,D,0.00554795034724531
222,E,0.00590293924064178
222,F,0.00590293924064178
222,G,0.00590293924064178
222,J,29556.6
222,K,1903.37626636961
222,L,425.943674590989
222,N,3223.56759423466
222,O,4610.85968628637
222,P,2414.22849555631
222,T,2127.619844211
222,Z,99.734
222,AF,1664.11988720118
222,AG,325.138226749943
223,D,0.00538361180457825
223,E,0.00582244292745138
223,F,0.00582244292745138
223,G,0.00582244292745138
223,J,29864
223,K,1926.45046033837
223,L,421.4434595
223,N,3251.6151609097
223,O,4642.55171794178
223,P,2438.81450200309
223,T,2134.55317754433
223,Z,92.712
223,AF,1680.00204124785
223,AG,326.586628889433
224,D,0.00532530567014544
224,E,0.00572449316892887
224,F,0.00572449316892887
224,G,0.00572449316892887
224,J,30165.4
224,K,1950.10624108514
224,L,424.146511694223
224,N,3278.64667947144
224,O,4674.48863793118
224,P,2462.44011672197
224,T,2141.48651087767
224,Z,93.712
224,AF,1696.06227169245
224,AG,328.041483265495
225,D,0.00533966582888934
225,E,0.00567091947478038
225,F,0.00567091947478038
225,G,0.00567091947478038
225,J,30451.4
225,K,1964.03127090479
225,L,422.990499243275
225,N,3305.65502820041
225,O,4706.67259288094
225,P,2485.75498407887
225,T,2132.180844211
225,Z,99.051
225,AF,1716.89000384285
225,AG,330.193673423702
226,D,0.00396062435642164
226,E,0.0056830615202379
226,F,0.0056830615202379
226,G,0.0056830615202379
226,J,30717.6
226,K,1976.60742450692
226,L,425.873223494629
226,N,3333.78857672717
226,O,4789.764743375
226,P,2506.61206147365
226,T,2176.11106217219
226,Z,97.901
226,AF,1738.00706068685
226,AG,332.359983510982
227,D,0
227,E,0
227,F,0
227,G,0
227,H,0
227,J,30987.6
228,J,31268.7
229,J,31552.8
230,J,31832.8
231,J,32110.2
232,J,32390.9
233,J,32675.3
234,J,32964
235,J,33259.6
236,J,33561.6
237,J,33869.9
238,J,34183.6
239,J,34502.8
240,J,34826.7
241,J,35154.4
242,J,35485.1
243,J,35819.3
244,J,36157.1
245,J,36498.4
246,J,36842.8
247,J,37190.4
248,J,37541.4
249,J,37894.8
250,J,38250.6
251,J,38609.4
252,J,38971.4
253,J,39336.2
254,J,39703.7
255,J,40074.1
256,J,40447.6
257,J,40824.3
258,J,41204.5
259,J,41588
260,J,41975
"@

$rows = $updates -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $rows) {
    $parts = $line.Trim() -split ","
    $rowNum = $parts[0]
    $colLetter = $parts[1]
    $newValue = [double]$parts[2]
    $ws.Range("$colLetter$rowNum").Value = $newValue
}
